$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Predicted DM/cm of twigs (column C, "Leaf weight") and branches (column D, "twigs weight")
# for each branch code, rows 2-23.
$values = @{
    2  = @(4.455, 1.028)
    3  = @(7.503, 2.485)
    4  = @(7.458, 2.29)
    5  = @(8.987, 3.168)
    6  = @(7.933, 2.544)
    7  = @(11.677, 8.728)
    8  = @(10.223, 4.96)
    9  = @(11.34, 4.887)
    10 = @(11.094, 4.614)
    11 = @(5.264, 1.152)
    12 = @(4.06, 0.811)
    13 = @(5.215, 0.92)
    14 = @(3.205, 1.124)
    15 = @(2.122, 0.924)
    16 = @(3.204, 1.282)
    17 = @(2.437, 0.775)
    18 = @(8.121, 2.395)
    19 = @(12.663, 3.765)
    20 = @(11.397, 12.006)
    21 = @(6.125, 2.648)
    22 = @(8.819, 3.733)
    23 = @(3.852, 1.447)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 3).Value = $pair[0]
    $ws.Cells.Item($row, 4).Value = $pair[1]
}

# Move the active selection to C12, matching the saved view state.
$ws.Range("C12").Select()
